$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 8414
$ws.Range("F8").Value = 2301
$ws.Range("F10").Value = 125
$ws.Range("F11").Value = 7670
$ws.Range("G11").Value = 78.2
$ws.Range("F12").Value = 7859
$ws.Range("F13").Value = 5033
$ws.Range("F15").Value = 16
$ws.Range("F16").Value = 629
$ws.Range("F17").Value = 5391
$ws.Range("F20").Value = 144
$ws.Range("F22").Value = 963
$ws.Range("F23").Value = 1500
$ws.Range("F24").Value = 2111
$ws.Range("F25").Value = 29
$ws.Range("F26").Value = 209
$ws.Range("F27").Value = 261
$ws.Range("F28").Value = 1103
$ws.Range("F32").Value = 798
$ws.Range("F33").Value = 1302
$ws.Range("F34").Value = 469
$ws.Range("F35").Value = 24
$ws.Range("F38").Value = 28
$ws.Range("F41").Value = 2510

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 7808
$ws.Range("F9").Value = 19
$ws.Range("F19").Value = 48
$ws.Range("F30").Value = 14
$ws.Range("F33").Value = 2
$ws.Range("F41").Value = 155

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 1644
$ws.Range("F7").Value = 694
$ws.Range("F8").Value = 2444
$ws.Range("F9").Value = 9492
$ws.Range("F10").Value = 1818
$ws.Range("F11").Value = 197
$ws.Range("F12").Value = 121
$ws.Range("F15").Value = 306
$ws.Range("F16").Value = 2611
$ws.Range("F17").Value = 291
$ws.Range("F18").Value = 107

$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 8414
$ws.Range("F6").Value = 694
$ws.Range("F7").Value = 1818
$ws.Range("F9").Value = 306
$ws.Range("F10").Value = 2611
$ws.Range("F11").Value = 291
$ws.Range("F12").Value = 5033
$ws.Range("F13").Value = 629
$ws.Range("F15").Value = 144
$ws.Range("F16").Value = 107
$ws.Range("F17").Value = 963
$ws.Range("F18").Value = 1500
$ws.Range("F22").Value = 19
$ws.Range("F26").Value = 261
$ws.Range("F30").Value = 798
$ws.Range("F32").Value = 1302
$ws.Range("F35").Value = 469
$ws.Range("F40").Value = 14
$ws.Range("F46").Value = 2510
